$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4434269841259209
$ws.Range("C2").Value = 0.04220382097970798
$ws.Range("D2").Value = 0.07802026236547022
$ws.Range("E2").Value = 0.4067577100716733
$ws.Range("G2").Value = 1.889026574525104
$ws.Range("H2").Value = 1.504283462932847
$ws.Range("K2").Value = 0.4086080624752526
$ws.Range("B3").Value = 0.407707959379735
$ws.Range("C3").Value = 0.03660660074201871
$ws.Range("D3").Value = 0.07082813792298737
$ws.Range("E3").Value = 0.3547810395071878
$ws.Range("G3").Value = 1.80454844583582
$ws.Range("H3").Value = 1.468146982374435
$ws.Range("K3").Value = 0.3698983548372894
$ws.Range("B4").Value = 0.3861192475179678
$ws.Range("C4").Value = 0.03317828547264412
$ws.Range("D4").Value = 0.06645224281474782
$ws.Range("E4").Value = 0.3229862608060898
$ws.Range("G4").Value = 1.753284275759285
$ws.Range("H4").Value = 1.446407518835201
$ws.Range("K4").Value = 0.3464372613667308
$ws.Range("B5").Value = 0.3774072398058763
$ws.Range("C5").Value = 0.03178312667584748
$ws.Range("D5").Value = 0.06467898304002517
$ws.Range("E5").Value = 0.3100571262850735
$ws.Range("G5").Value = 1.732543573576919
$ws.Range("H5").Value = 1.437660231697464
$ws.Range("K5").Value = 0.3369529706430683
$ws.Range("B6").Value = 0.3759657666438727
$ws.Range("C6").Value = 0.03155157089616978
$ws.Range("D6").Value = 0.06438513198555995
$ws.Range("E6").Value = 0.3079118428106256
$ws.Range("G6").Value = 1.729108575748853
$ws.Range("H6").Value = 1.436214468371048
$ws.Range("K6").Value = 0.3353826947219716
$ws.Range("B7").Value = 0.3860014087395882
$ws.Range("C7").Value = 0.03315946241660583
$ws.Range("D7").Value = 0.06642828788305621
$ws.Range("E7").Value = 0.3228117858470512
$ws.Range("G7").Value = 1.753003955849096
$ws.Range("H7").Value = 1.446289098922165
$ws.Range("K7").Value = 0.346309045110786
$ws.Range("B8").Value = 0.431039478636734
$ws.Range("C8").Value = 0.04027203022099002
$ws.Range("D8").Value = 0.07553201439129964
$ws.Range("E8").Value = 0.3888099342955087
$ws.Range("G8").Value = 1.859771285243227
$ws.Range("H8").Value = 1.491729873330257
$ws.Range("K8").Value = 0.3951967403294816
$ws.Range("B9").Value = 0.5221153435555834
$ws.Range("C9").Value = 0.05429638981948415
$ws.Range("D9").Value = 0.09370943243905572
$ws.Range("E9").Value = 0.5192945316018864
$ws.Range("G9").Value = 2.074070071428963
$ws.Range("H9").Value = 1.584450513604139
$ws.Range("K9").Value = 0.4935444995633986
$ws.Range("B10").Value = 0.590766457380937
$ws.Range("C10").Value = 0.06466169599059413
$ws.Range("D10").Value = 0.107273510326408
$ws.Range("E10").Value = 0.6159814252814328
$ws.Range("G10").Value = 2.234705258035859
$ws.Range("H10").Value = 1.654855060745831
$ws.Range("K10").Value = 0.5673823253993078
$ws.Range("B11").Value = 0.6223880691990473
$ws.Range("C11").Value = 0.06939384698712558
$ws.Range("D11").Value = 0.1134920562779627
$ws.Range("E11").Value = 0.6601830527532684
$ws.Range("G11").Value = 2.308516274226008
$ws.Range("H11").Value = 1.687397692495495
$ws.Range("K11").Value = 0.6013325148822446
$ws.Range("B12").Value = 0.6344196348586877
$ws.Range("C12").Value = 0.07118847818139784
$ws.Range("D12").Value = 0.1158539607909006
$ws.Range("E12").Value = 0.6769554281978003
$ws.Range("G12").Value = 2.336575643128924
$ws.Range("H12").Value = 1.699796161453037
$ws.Range("K12").Value = 0.6142416418597065
$ws.Range("B13").Value = 0.6318258673706225
$ws.Range("C13").Value = 0.07080184977991166
$ws.Range("D13").Value = 0.1153449652963303
$ws.Range("E13").Value = 0.6733416320247869
$ws.Range("G13").Value = 2.330527676118322
$ws.Range("H13").Value = 1.697122562106131
$ws.Range("K13").Value = 0.6114590600020335
$ws.Range("B14").Value = 0.6233767633909792
$ws.Range("C14").Value = 0.06954143753617359
$ws.Range("D14").Value = 0.1136862289103959
$ws.Range("E14").Value = 0.6615622250334923
$ws.Range("G14").Value = 2.310822540400352
$ws.Range("H14").Value = 1.688416206403019
$ws.Range("K14").Value = 0.6023934898821324
$ws.Range("B15").Value = 0.6182089095203764
$ws.Range("C15").Value = 0.06876975318820655
$ws.Range("D15").Value = 0.1126711309400861
$ws.Range("E15").Value = 0.6543515328023091
$ws.Range("G15").Value = 2.298766818134141
$ws.Range("H15").Value = 1.683093151426647
$ws.Range("K15").Value = 0.5968474896675673
$ws.Range("B16").Value = 0.5887078638275511
$ws.Range("C16").Value = 0.06435279928577131
$ws.Range("D16").Value = 0.1068680958023407
$ws.Range("E16").Value = 0.613097374723452
$ws.Range("G16").Value = 2.229896602081993
$ws.Range("H16").Value = 1.652738792596324
$ws.Range("K16").Value = 0.5651709511374747
$ws.Range("B17").Value = 0.5707108567624459
$ws.Range("C17").Value = 0.06164764093698238
$ws.Range("D17").Value = 0.1033205659773415
$ws.Range("E17").Value = 0.5878471614115313
$ws.Range("G17").Value = 2.187837712572446
$ws.Range("H17").Value = 1.634250185104406
$ws.Range("K17").Value = 0.5458315955117712
$ws.Range("B18").Value = 0.5603962906487254
$ws.Range("C18").Value = 0.06009328378560497
$ws.Range("D18").Value = 0.1012846464611528
$ws.Range("E18").Value = 0.5733442819012424
$ws.Range("G18").Value = 2.163715781931018
$ws.Range("H18").Value = 1.623664452895241
$ws.Range("K18").Value = 0.5347420257095337
$ws.Range("B19").Value = 0.5569102585387213
$ws.Range("C19").Value = 0.05956726916100763
$ws.Range("D19").Value = 0.1005960919398916
$ws.Range("E19").Value = 0.5684372713150054
$ws.Range("G19").Value = 2.155560318634286
$ws.Range("H19").Value = 1.620088587072843
$ws.Range("K19").Value = 0.5309930824373623
$ws.Range("B20").Value = 0.5726228507632811
$ws.Range("C20").Value = 0.06193544461264366
$ws.Range("D20").Value = 0.1036977371113466
$ws.Range("E20").Value = 0.5905329608268204
$ws.Range("G20").Value = 2.192307768490537
$ws.Range("H20").Value = 1.636213308527886
$ws.Range("K20").Value = 0.5478867854583882
$ws.Range("B21").Value = 0.625856912623334
$ws.Range("C21").Value = 0.06991157685241944
$ws.Range("D21").Value = 0.1141732467583978
$ws.Range("E21").Value = 0.66502117233901
$ws.Range("G21").Value = 2.31660744342912
$ws.Range("H21").Value = 1.690971421476149
$ws.Range("K21").Value = 0.6050548229328001
$ws.Range("B22").Value = 0.6609819332511506
$ws.Range("C22").Value = 0.07514010557449069
$ws.Range("D22").Value = 0.1210609135610952
$ws.Range("E22").Value = 0.7139040714992007
$ws.Range("G22").Value = 2.398479566465994
$ws.Range("H22").Value = 1.727198369311054
$ws.Range("K22").Value = 0.6427264284079399
$ws.Range("B23").Value = 0.642204266441297
$ws.Range("C23").Value = 0.07234802920930861
$ws.Range("D23").Value = 0.1173810061546021
$ws.Range("E23").Value = 0.6877950638050692
$ws.Range("G23").Value = 2.354723854353722
$ws.Range("H23").Value = 1.707822767676021
$ws.Range("K23").Value = 0.6225917689632467
$ws.Range("B24").Value = 0.5717583381856457
$ws.Range("C24").Value = 0.06180532586647303
$ws.Range("D24").Value = 0.1035272068327799
$ws.Range("E24").Value = 0.5893186681552436
$ws.Range("G24").Value = 2.190286674868645
$ws.Range("H24").Value = 1.635325644647594
$ws.Range("K24").Value = 0.5469575441968288
$ws.Range("B25").Value = 0.4971759666509286
$ws.Range("C25").Value = 0.0504927001627351
$ws.Range("D25").Value = 0.08875603608659333
$ws.Range("E25").Value = 0.4838636042419608
$ws.Range("G25").Value = 2.015549892934047
$ws.Range("H25").Value = 1.558972166355062
$ws.Range("K25").Value = 0.4666662221212619
